$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.949.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "'1.746.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'232.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.5197"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").Value = "'0.2814"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.22%  "

$ws.Range("D9").Value = "'39.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.09%  "

$ws.Range("D10").Value = "'0.06129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").Value = "'1.749.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").Value = "'0.07042"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").Value = "'15.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").Value = "'0.6452"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").Value = "'4.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "'77.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").Value = "'25.914.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").Value = "'11.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").Value = "'0.000006608"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.56%  "

$ws.Range("D22").Value = "'1.977.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'4.141"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("D24").Value = "'8.653"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.74%  "

$ws.Range("D25").Value = "'5.157"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("D26").Value = "'139.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "

$ws.Range("E27").Value = "  +3.46%  "

$ws.Range("D28").Value = "'15.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").Value = "'1.821"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").Value = "'102.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").Value = "'0.08270"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").Value = "'3.674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").Value = "'3.440"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").Value = "'0.04480"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").Value = "'2.611"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").Value = "'0.9893"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("D37").Value = "'0.6185"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.06%  "

$ws.Range("D38").Value = "'2.669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("D39").Value = "'0.01590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.67%  "

$ws.Range("D40").Value = "'1.924"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").Value = "'1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").Value = "'100.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").Value = "'0.3850"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.062"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.95%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.7289"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.33%  "

$ws.Range("D46").Value = "'0.05449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("D47").Value = "'6.329"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.48%  "

$ws.Range("D48").Value = "'0.1124"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("D49").Value = "'53.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").Value = "'30.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("E51").Value = "  +1.88%  "
